$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row "101, 1, hombre" (worksheet row 2) was removed and the rows below it
# shifted up. Deleting the entire row achieves exactly that shift.
$ws.Rows.Item(2).Delete()

# Update the active selection as recorded in the saved file.
$ws.Range("E10").Select()
